$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 already has cells A2 and B2 with a text number format (style index 1).
# Reuse that same format for the other new data cells that need it (A2, B2, C2),
# matching the original author's text-formatted ID / ServerID / Name columns.
$textFormat = $ws.Range("A2").NumberFormat

# F2 = IP column -> "127.0.0.1" (plain text, default format). Filled in first so
# new shared-string entries land in the same order as the authored workbook.
$ws.Range("F2").Value = "127.0.0.1"

# A2 = ID column -> "ProxyServer_1"
$ws.Range("A2").Value = "ProxyServer_1"

# B2 = ServerID column -> "000105001" (kept as text because of the existing numFmt)
$ws.Range("B2").Value = "000105001"

# C2 = Name column -> "ProxyServer_1" (new cell, needs same text format as A2/B2)
$ws.Range("C2").NumberFormat = $textFormat
$ws.Range("C2").Value = "ProxyServer_1"

# D2 = MaxOnline column -> 5000 (plain number, default format)
$ws.Range("D2").Value = 5000

# E2 = CpuCount column -> 1 (plain number, default format)
$ws.Range("E2").Value = 1

# G2 = Port column -> 5001 (plain number, default format)
$ws.Range("G2").Value = 5001

# Update the active selection to G4, as recorded in the saved view state.
$ws.Range("G4").Select() | Out-Null
